$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01068002311887565
$ws.Range("C2").Value = 0.008480156111028611
$ws.Range("D2").Value = 0.009175445401755538
$ws.Range("E2").Value = 0.009419315002041663
$ws.Range("B3").Value = 5.523028350170109
$ws.Range("C3").Value = 5.949757187784542
$ws.Range("D3").Value = 7.409635965204234
$ws.Range("E3").Value = 7.923920638043707
$ws.Range("B4").Value = -0.004754311418771289
$ws.Range("C4").Value = -0.005389076851612592
$ws.Range("D4").Value = -0.00585536679281113
$ws.Range("E4").Value = -0.005670665460181294
$ws.Range("B5").Value = -2.95150647669481
$ws.Range("C5").Value = -4.754240764293362
$ws.Range("D5").Value = -5.585856170897704
$ws.Range("E5").Value = -5.992451452608591
$ws.Range("B6").Value = 0.00117625630828122
$ws.Range("C6").Value = 0.0009029167040204621
$ws.Range("D6").Value = 0.001258063373843503
$ws.Range("B7").Value = 1.854481496400236
$ws.Range("C7").Value = 0.9257860547129588
$ws.Range("D7").Value = 0.9101392953611729
$ws.Range("B8").Value = 0.008705333759653674
$ws.Range("C8").Value = 0.008192026803880783
$ws.Range("D8").Value = 0.008836394527485313
$ws.Range("E8").Value = 0.009776342004670801
$ws.Range("B9").Value = 4.547025827400423
$ws.Range("C9").Value = 5.824689285748136
$ws.Range("D9").Value = 7.041076972100504
$ws.Range("E9").Value = 7.803803585805547
$ws.Range("B10").Value = -0.006451957805491795
$ws.Range("C10").Value = -0.006356541366461522
$ws.Range("D10").Value = -0.006696391702962223
$ws.Range("E10").Value = -0.005730842260671652
$ws.Range("B11").Value = -3.858141706739029
$ws.Range("C11").Value = -5.400859541094776
$ws.Range("D11").Value = -6.478674860957485
$ws.Range("E11").Value = -6.379229698710517
$ws.Range("B12").Value = 0.0002720136651238718
$ws.Range("C12").Value = 0.0002409150076297867
$ws.Range("D12").Value = 0.0004160474054679415
$ws.Range("B13").Value = 0.4322938521137515
$ws.Range("C13").Value = 0.2413693795401204
$ws.Range("D13").Value = 0.3057493790834024
$ws.Range("B14").Value = 0.008557938219821676
$ws.Range("C14").Value = 0.00761211709070698
$ws.Range("D14").Value = 0.008439940881558992
$ws.Range("E14").Value = 0.009506509389972186
$ws.Range("B15").Value = 4.336635900314954
$ws.Range("C15").Value = 5.466529391994337
$ws.Range("D15").Value = 6.705800154687222
$ws.Range("E15").Value = 7.462273818520735
$ws.Range("B16").Value = -0.006695666715188392
$ws.Range("C16").Value = -0.006681371175403872
$ws.Range("D16").Value = -0.006641293356938142
$ws.Range("E16").Value = -0.00569243078486502
$ws.Range("B17").Value = -3.744964528610712
$ws.Range("C17").Value = -5.432774179440074
$ws.Range("D17").Value = -6.37000633007165
$ws.Range("E17").Value = -6.026796462879193
$ws.Range("B18").Value = 0.0001140004600902067
$ws.Range("C18").Value = -0.0002100742774031683
$ws.Range("D18").Value = 0.0002437057523812515
$ws.Range("B19").Value = 0.1666436678023185
$ws.Range("C19").Value = -0.2157019619676573
$ws.Range("D19").Value = 0.184412967581387
$ws.Range("B20").Value = 0.008694176507028759
$ws.Range("C20").Value = 0.008243191096786459
$ws.Range("D20").Value = 0.009032967302944159
$ws.Range("E20").Value = 0.009611170856552198
$ws.Range("B21").Value = 4.42452692558305
$ws.Range("C21").Value = 5.739419124279273
$ws.Range("D21").Value = 7.100177912027119
$ws.Range("E21").Value = 8.003958951874534
$ws.Range("B22").Value = -0.006366746531862169
$ws.Range("C22").Value = -0.006750748498099018
$ws.Range("D22").Value = -0.006539269801902266
$ws.Range("E22").Value = -0.005772245651577344
$ws.Range("B23").Value = -3.617560667053156
$ws.Range("C23").Value = -5.44956356670746
$ws.Range("D23").Value = -5.984021443800612
$ws.Range("E23").Value = -5.817743972605122
$ws.Range("B24").Value = 0.0002525253575056304
$ws.Range("C24").Value = 0.00007639543261224925
$ws.Range("D24").Value = 0.0004237641230775907
$ws.Range("B25").Value = 0.3789261319459177
$ws.Range("C25").Value = 0.0760645673010225
$ws.Range("D25").Value = 0.2832520318358187
